$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ColorScheme
Write-Output $cs
Write-Output $cs.Count
